# Delete unused variables in Molten Salt Tower Parasitics UI page
# Adds 7 new "Deleted variable" rows (36-42) to the "SAM Variable Changes"
# sheet, mirroring the format of the row directly above (row 35).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SAM Variable Changes")
$ws.Activate()

# New deleted-variable names (column C) and their "reason deleted" text
# (column F). All rows share the same Input Page (column E = "Molten Salt
# Tower Parasitics"), Type (column A = "Deleted variable"), Variable Type
# (column B = "number"), Version-Upgrader flag (column G = "N/A") and
# Ty marker (column H = "Ty").
#
# UseShadedReason marks rows whose "reason deleted" cell should use the
# shaded style (matching other "not used" rows elsewhere in the sheet)
# instead of the plain style inherited from row 35.
$newVars = @(
    @{ Name = "P_storage_pump";    Reason = "storage HTF = rec/pc HTF (no storage HX), so no pumping losses"; UseShadedReason = $false },
    @{ Name = "storage_bypass";    Reason = "storage HTF = rec/pc HTF (no storage HX), so no pumping losses"; UseShadedReason = $false },
    @{ Name = "recirc_source";     Reason = "not used"; UseShadedReason = $true },
    @{ Name = "recirc_htf_eff";    Reason = "not used"; UseShadedReason = $true },
    @{ Name = "flow_from_storage"; Reason = "not used"; UseShadedReason = $true },
    @{ Name = "P_hot_tank";        Reason = "not used"; UseShadedReason = $true },
    @{ Name = "csp.pt.par.bop_c1"; Reason = "not used"; UseShadedReason = $true }
)

$startRow = 36

for ($i = 0; $i -lt $newVars.Count; $i++) {
    $row = $startRow + $i
    $prevRow = $row - 1

    # Duplicate the row immediately above (keeps per-cell styles identical
    # to the existing "Deleted variable" rows) and shift everything below
    # down by one.
    $ws.Rows($prevRow.ToString() + ":" + $prevRow.ToString()).Copy()
    $ws.Rows($row.ToString() + ":" + $row.ToString()).Insert(-4121)

    $entry = $newVars[$i]

    $ws.Range("C" + $row).Value = $entry.Name

    if ($entry.UseShadedReason) {
        # Rows whose reason is "not used" pick up the shaded fill style
        # (same as the other shaded cells in this column) rather than the
        # plain style copied from row 35.
        $ws.Range("F34").Copy()
        $ws.Range("F" + $row).PasteSpecial(-4122)
    }

    $ws.Range("F" + $row).Value = $entry.Reason
}

# Update the saved view state to match where the user ended up after
# entering the new rows.
$ws.Range("A43").Select()
